# Add files via upload
#
# Recreates the target state of "Шеймин 4.xlsx": three new localization
# rows (4, 5, 6) are appended below the existing row 3, plus a trailing
# row (7) that only carries a filename in column A. Row 3 itself is
# extended with a filename in column A and gets a thin bottom border
# under its C:E cells (since it is no longer the last/bottom row of the
# table). New shared strings are entered column-by-column (C, then D,
# then E, then A) which is the order a translator would paste batches of
# text into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (English) ------------------------------------------------
$ws.Range("C4").Value = ' Climbing a mountain is\nhard work…'
$ws.Range("C5").Value = ' When you get to the summit,\nthe sense of achievement will leave you\nspeechless, don\''t you think?'

# --- Column D (Russian translation) ------------------------------------
$ws.Range("D4").Value = ' Восхождение на гору даётся\nнепросто...'
$ws.Range("D5").Value = ' Но когда ты доходишь до вершины\nгоры, тебя до немоты переполняет чувство\nуспеха, верно?'

# --- Column E (converted/encoded string) --------------------------------
$ws.Range("E4").Value = ' Âïòöïçäåîéå îà ãïñô äàæóòÿ\nîåðñïòóï...'
$ws.Range("E5").Value = ' Îï ëïãäà óú äïöïäéšû äï âåñšéîú\nãïñú, óåáÿ äï îåíïóú ðåñåðïìîÿåó œôâòóâï\nôòðåöà, âåñîï?'

# --- Column A (script filenames), top to bottom --------------------------
$ws.Range("A3").Value = "SCRIPT/D73P11A/us0205.ssb"
$ws.Range("A4").Value = "SCRIPT/D73P11A/us0304.ssb"
$ws.Range("A5").Value = "SCRIPT/D73P11A/us0405.ssb"
$ws.Range("A6").Value = "SCRIPT/D73P11A/us2004.ssb"
$ws.Range("A7").Value = "SCRIPT/D73P11A/us2104.ssb"

# --- Column B (line numbers) ---------------------------------------------
$ws.Range("B3").Value = 44
$ws.Range("B4").Value = 18
$ws.Range("B5").Value = 21

# --- Row 3 is no longer the final row, so it gains a thin divider
#     border under C3:E3 -------------------------------------------------
$ws.Range("C3:E3").Borders.Item(9).LineStyle = 1
$ws.Range("C3:E3").Borders.Item(9).Weight = 2

# --- Row heights: rows 3-7 all grow to 43.2pt -----------------------------
$ws.Range("A3:E7").RowHeight = 43.2

# --- Final selection, matching the saved view state -----------------------
$ws.Range("E5").Select()
